$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above "Fiscal Year" (row 3) to hold the new "Report Type"
# field. This pushes the existing field-label rows and the column-header row
# down by one (the header row moves from row 19 to row 20).
$ws.Rows("3:3").Insert()

# The newly inserted row inherits formatting from the row above it (the
# title row), so copy the formatting of the sibling field-label cell
# (now at A4, "Fiscal Year") onto the new A3 cell, then set its text.
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Report Type"

# Update the cell that is selected/shown as active when the sheet is opened.
$ws.Range("C8").Select() | Out-Null
